$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8449057340621948
$ws.Range("B1").Value = 1.298365592956543
$ws.Range("C1").Value = 5.194911956787109
$ws.Range("D1").Value = 1.622901678085327
$ws.Range("E1").Value = 0.9331645965576172
